# Refresh the cryptocurrency price/volume table on Sheet1 (rows 2-51,
# columns B-E) with the latest values from this run's scrape, per the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Columns D (Price) and E (Volume(1h)) are stored as plain text, not real
# numbers -- e.g. "43.691.03" (dotted thousands) or "  -0.38%  " (padded
# percent string), and some Price values look like ordinary decimals
# (e.g. "71.86", "0.581", "174.90"). A plain Range.Value assignment of a
# decimal-looking string gets auto-coerced by Excel into a Double (so
# "174.90" silently becomes 174.9, losing the trailing zero / text type).
# To keep those as literal text we write them with a leading apostrophe
# (forces text / quote-prefix entry) and then immediately ClearFormats()
# on just that cell so the transient "quote prefix" cell style Excel
# applies doesn't linger as a spurious formatting change.
#
# Rows 38/39 are a full content swap: THORChain now outranks LidoDAOToken,
# so columns B/C/D/E of both rows are rewritten in place (column A, the
# rank index, is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.691.03"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.335.10"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("E6").Value = "  -3.91%  "
$ws.Range("D7").Value = "'71.86"
$ws.Range("E7").Value = "  -5.55%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.581"
$ws.Range("E9").Value = "  -7.08%  "
$ws.Range("E10").Value = "  -4.13%  "
$ws.Range("D11").Value = "'58.30"
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "'7.09"
$ws.Range("E14").Value = "  -5.54%  "
$ws.Range("D15").Value = "2.682.23"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("E16").Value = "  -4.96%  "
$ws.Range("D17").Value = "'0.891"
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("D18").Value = "2.333.61"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").Value = "43.613.40"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("E20").Value = "  -3.67%  "
$ws.Range("D21").Value = "'77.66"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "'6.58"
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").Value = "'250.19"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("E24").Value = "  +7.57%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  +2.57%  "
$ws.Range("D27").Value = "'2.48"
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("D28").Value = "'10.28"
$ws.Range("E28").Value = "  -8.31%  "
$ws.Range("D29").Value = "'2.26"
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("D30").Value = "'174.90"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "'22.05"
$ws.Range("E31").Value = "  -4.84%  "
$ws.Range("D32").Value = "'0.126"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("D34").Value = "'0.0731"
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("E35").Value = "  -4.84%  "
$ws.Range("D36").Value = "'5.34"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").Value = "'6.34"
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'2.36"
$ws.Range("E39").Value = "  -3.33%  "
$ws.Range("D40").Value = "'5.52"
$ws.Range("E40").Value = "  +22.25%  "
$ws.Range("D41").Value = "'0.0268"
$ws.Range("E41").Value = "  -2.94%  "
$ws.Range("D42").Value = "'64.68"
$ws.Range("E42").Value = "  +18.98%  "
$ws.Range("D43").Value = "'9.17"
$ws.Range("E43").Value = "  +3.10%  "
$ws.Range("E44").Value = "  +3.67%  "
$ws.Range("D45").Value = "'18.73"
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("E46").Value = "  -2.94%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  -2.88%  "
$ws.Range("E49").Value = "  -3.46%  "
$ws.Range("E50").Value = "  +3.94%  "
$ws.Range("D51").Value = "'97.43"
$ws.Range("E51").Value = "  -4.15%  "

# Strip the transient quote-prefix formatting picked up by the apostrophe
# trick above, cell by cell (a combined multi-area Range(...) string was
# unreliable for ClearFormats in this host, so each cell is cleared
# individually to be safe).
$ws.Range("D7").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D51").ClearFormats()
